# The "Ano" column (B) held date-serial numbers formatted with a custom
# date/time number format. Convert those cells (rows 2-37) to plain text
# values written as "dd/mm/yyyy", dropping the custom date formatting so
# the cells fall back to the default (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $serial = $cell.Value2
    $dt = [DateTime]::FromOADate($serial)
    $text = $dt.ToString("dd/MM/yyyy")

    $cell.Value = "'" + $text
    $cell.ClearFormats()
}
